$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch every target cell with a text number-format so the new value is
# stored as text (matching the original shared-string cell type), then
# clear the format back to General so no style is left behind.
$targetRange = $ws.Range("C2,D2,E2,R2,S2,T2,U2,C3,D3,E3,R3,S3,T3,U3,C4,D4,E4,R4,S4,T4,U4,C5,D5,E5,R5,S5,T5,U5,C6,D6,E6,R6,S6,T6,U6,C7,D7,E7,R7,S7,T7,C8,D8,E8,R8,T8,U8,V8,C9,D9,E9,R9,S9,T9,U9,V9,W9,C10,D10,E10,R10,S10,T10,U10,V10,C11,D11,E11,R11,S11,T11,U11")
$targetRange.NumberFormat = "@"

$ws.Range("C2").Value = "-6.911871570822404"
$ws.Range("D2").Value = "-0.21521026486142603"
$ws.Range("E2").Value = "52.98480893488265"
$ws.Range("R2").Value = "-75.28451095800045"
$ws.Range("S2").Value = "-41.70211347736648"
$ws.Range("T2").Value = "41.27169294764363"
$ws.Range("U2").Value = "48.06744520443369"
$ws.Range("C3").Value = "-0.3543272665957655"
$ws.Range("D3").Value = "1.1170705284749474"
$ws.Range("E3").Value = "42.304583654301275"
$ws.Range("R3").Value = "-45.85820751037421"
$ws.Range("S3").Value = "41.56662012781063"
$ws.Range("T3").Value = "42.20675738704125"
$ws.Range("U3").Value = "-39.332479070860735"
$ws.Range("C4").Value = "33.92451846423617"
$ws.Range("D4").Value = "51.52438930839696"
$ws.Range("E4").Value = "44.88061590021839"
$ws.Range("R4").Value = "74.2225939885502"
$ws.Range("S4").Value = "-41.57329874839941"
$ws.Range("T4").Value = "59.064698719608465"
$ws.Range("U4").Value = "43.98407989718545"
$ws.Range("C5").Value = "22.181912640697156"
$ws.Range("D5").Value = "46.96518492533878"
$ws.Range("E5").Value = "48.05680886903047"
$ws.Range("R5").Value = "55.60490802532658"
$ws.Range("S5").Value = "-60.80762731321552"
$ws.Range("T5").Value = "45.37679231239601"
$ws.Range("U5").Value = "48.55357753828155"
$ws.Range("C6").Value = "-7.045877421782528"
$ws.Range("D6").Value = "-4.893257402698325"
$ws.Range("E6").Value = "46.52843976698257"
$ws.Range("R6").Value = "-63.420982772657304"
$ws.Range("S6").Value = "45.02398789092384"
$ws.Range("T6").Value = "-42.0476104834624"
$ws.Range("U6").Value = "32.26109567806575"
$ws.Range("C7").Value = "-35.69223966761786"
$ws.Range("D7").Value = "-64.12503882579095"
$ws.Range("E7").Value = "50.72933203675473"
$ws.Range("R7").Value = "-78.5186279093085"
$ws.Range("S7").Value = "35.56694773224588"
$ws.Range("T7").Value = "-64.12503882579095"
$ws.Range("C8").Value = "6.070738168030725"
$ws.Range("D8").Value = "39.363418170605705"
$ws.Range("E8").Value = "44.02958802570318"
$ws.Range("R8").Value = "-50.014966639368545"
$ws.Range("T8").Value = "39.363418170605705"
$ws.Range("U8").Value = "44.842795260311505"
$ws.Range("V8").Value = "41.71651241797505"
$ws.Range("C9").Value = "-13.351267345975979"
$ws.Range("D9").Value = "-12.371830583187212"
$ws.Range("E9").Value = "48.1214049460452"
$ws.Range("R9").Value = "30.720433098962136"
$ws.Range("S9").Value = "-71.62470029295248"
$ws.Range("T9").Value = "-52.75348378699806"
$ws.Range("U9").Value = "28.009822620623638"
$ws.Range("V9").Value = "-58.306700010311275"
$ws.Range("W9").Value = "43.84702429482017"
$ws.Range("C10").Value = "-24.21582632878154"
$ws.Range("D10").Value = "-47.6048969763473"
$ws.Range("E10").Value = "51.567206126551696"
$ws.Range("R10").Value = "-92.16758154716233"
$ws.Range("S10").Value = "-53.43252148827322"
$ws.Range("T10").Value = "38.378194594089344"
$ws.Range("U10").Value = "-47.6048969763473"
$ws.Range("V10").Value = "33.74767377378584"
$ws.Range("C11").Value = "0.5436438980255804"
$ws.Range("D11").Value = "1.1107880255880183"
$ws.Range("E11").Value = "39.31185151147227"
$ws.Range("R11").Value = "-41.49590915749956"
$ws.Range("S11").Value = "38.12753195734236"
$ws.Range("T11").Value = "41.448908698425846"
$ws.Range("U11").Value = "-35.90595590616633"

$targetRange.ClearFormats()

